$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A new run (Run 50) was appended as a new last column (BA), shifting the
# previous "Mean" column from AZ to BA and pushing new run-50 results into
# the former Mean column (AZ), along with the recalculated mean in BA.

# Copy AZ1's formatting (bold font, thin border, centered alignment) into
# the new BA1 header cell before overwriting the header text values.
$ws.Range("AZ1").Copy($ws.Range("BA1"))

$ws.Range("AZ1").Value = "Run 50"
$ws.Range("BA1").Value = "Mean"

# Fill in the Run 50 results (now in column AZ) and the recalculated Mean
# (now in the new column BA) for every data row.
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 52).Value = 29.08910829
    $ws.Cells.Item($r, 53).Value = 23.80975355
}
